$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Fecha 2022-06-02 week) values
$ws.Range("D2").Value = 44714
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 9000
$ws.Range("M2").Value = 9500
$ws.Range("P2").Value = 528

# Row 3 (Fecha 2022-08-31 week) values
$ws.Range("D3").Value = 44804
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 9500
$ws.Range("M3").Value = 9750
$ws.Range("P3").Value = 542
